$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIMULADOR_FII")

# Fix the PERFIL cell: it was mistakenly set to "WERER" (a typo, and not one of
# the valid dropdown choices CONSERVADOR / MODERADO / AGRESSIVO backed by the
# TAB_AUX!$D$3:$F$3 list). Set it to "CONSERVADOR" so the IFS/XLOOKUP driven
# table (Tabela2) and the pie chart stop evaluating to #N/A.
$ws.Range("C12").Value = "CONSERVADOR"

# Update the active selection to match the cell that was edited.
$ws.Range("C12").Select()
